$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain text values in the source data
# (t="inlineStr"), even though many of them look numeric (e.g. "0.9998").
# Temporarily force a Text number format while assigning the new values so
# Excel keeps storing them as strings instead of silently converting them
# to floating point numbers, then restore the original (default) format
# and cell style so no stray style attributes are left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.780.26"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "1.860.10"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "245.08"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "0.6422"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "47.57"
$ws.Range("E8").Value = "  +4.22%  "
$ws.Range("D9").Value = "0.07542"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").Value = "0.2976"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +6.42%  "
$ws.Range("D12").Value = "0.07676"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "1.869.32"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "0.6933"
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").Value = "84.04"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "0.000009843"
$ws.Range("E17").Value = "  +9.70%  "
$ws.Range("D18").Value = "6.132"
$ws.Range("E18").Value = "  +5.01%  "
$ws.Range("D19").Value = "29.785.27"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "2.113.71"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").Value = "237.02"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "1.0000"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "7.531"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "158.92"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "0.1425"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").Value = "8.564"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").Value = "0.06236"
$ws.Range("E30").Value = "  +6.47%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "1.296"
$ws.Range("E32").Value = "  +5.96%  "
$ws.Range("D33").Value = "4.157"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").Value = "4.113"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "1.904"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").Value = "1.175"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("D37").Value = "0.7304"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "2.610"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "2.825"
$ws.Range("D40").Value = "0.01787"
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("D41").Value = "1.212.84"
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "0.9237"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("D43").Value = "6.298"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "2.027.15"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").Value = "102.14"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "66.91"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000119"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.232"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.4060"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.670"
$ws.Range("E51").Value = "  +5.40%  "

# Restore default formatting/style on the Price column.
$ws.Range("D2:D51").Style = "Normal"
